# Applies the "Updated cryptos list" price/volume refresh to Sheet1.
# Column D ("Price") and E ("Volume(1h)") are stored as plain text in the
# sheet (values like "42.637.72" / "0.0000106" / "  +0.67%  "), so writes that
# look like numbers must be forced back to Text or Excel auto-converts them
# (e.g. "1.00" -> 1, "0.0000106" -> 1.06E-5), which would change both the
# stored value and the cell type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    # Force text storage for values that Excel would otherwise parse as a number.
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    # Drop the temporary Text number-format again so the cell keeps the same
    # (unstyled) look it had before - only the content should change.
    $r.ClearFormats()
}

function Set-PlainValue($cellRef, $value) {
    # Safe as-is: not numeric-looking, so Excel keeps it as text.
    $ws.Range($cellRef).Value = $value
}

Set-PlainValue 'D2' '42.488.64'
Set-PlainValue 'E2' '  +0.51%  '
Set-PlainValue 'D3' '2.301.48'
Set-PlainValue 'E3' '  +0.39%  '
Set-PlainValue 'E4' '  -0.29%  '
Set-TextValue 'D5' '316.38'
Set-PlainValue 'E5' '  +0.10%  '
Set-TextValue 'D6' '103.62'
Set-PlainValue 'E6' '  -0.51%  '
Set-TextValue 'D7' '0.629'
Set-PlainValue 'E7' '  -0.35%  '
Set-PlainValue 'E8' '  -0.03%  '
Set-PlainValue 'E9' '  -0.22%  '
Set-TextValue 'D10' '39.80'
Set-PlainValue 'E10' '  +0.63%  '
Set-PlainValue 'E11' '  -0.33%  '
Set-TextValue 'D12' '8.53'
Set-PlainValue 'E12' '  +2.88%  '
Set-PlainValue 'E13' '  +0.45%  '
Set-TextValue 'D14' '0.993'
Set-PlainValue 'E14' '  +3.36%  '
Set-TextValue 'D15' '15.36'
Set-PlainValue 'E15' '  +0.63%  '
Set-PlainValue 'D16' '2.650.83'
Set-PlainValue 'E16' '  +0.41%  '
Set-PlainValue 'D17' '2.295.65'
Set-PlainValue 'E17' '  -0.30%  '
Set-PlainValue 'D18' '42.451.94'
Set-PlainValue 'E18' '  +0.47%  '
Set-PlainValue 'E19' '  +3.88%  '
Set-PlainValue 'B20' 'InternetComputer(DFINITY)'
Set-PlainValue 'C20' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D20' '13.85'
Set-PlainValue 'E20' '  +33.03%  '
Set-PlainValue 'B21' 'ShibaInu'
Set-PlainValue 'C21' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D21' '0.0000106'
Set-PlainValue 'E21' '  +0.73%  '
Set-TextValue 'D22' '73.97'
Set-PlainValue 'E22' '  +0.90%  '
Set-PlainValue 'E23' '  -1.61%  '
Set-TextValue 'D24' '267.82'
Set-PlainValue 'E24' '  -3.42%  '
Set-PlainValue 'E25' '  -1.20%  '
Set-TextValue 'D26' '1.00'
Set-PlainValue 'E26' '  -0.29%  '
Set-PlainValue 'E27' '  +0.78%  '
Set-PlainValue 'E28' '  -3.82%  '
Set-TextValue 'D29' '22.67'
Set-PlainValue 'E29' '  -0.44%  '
Set-PlainValue 'E30' '  +14.26%  '
Set-TextValue 'D31' '37.69'
Set-PlainValue 'E31' '  +4.04%  '
Set-TextValue 'D32' '165.42'
Set-PlainValue 'E32' '  +0.96%  '
Set-PlainValue 'E33' '  +1.34%  '
Set-PlainValue 'E34' '  -5.67%  '
Set-TextValue 'D35' '0.132'
Set-PlainValue 'E35' '  -3.38%  '
Set-PlainValue 'E36' '  +0.67%  '
Set-PlainValue 'E37' '  +1.50%  '
Set-PlainValue 'E38' '  +1.68%  '
Set-TextValue 'D39' '3.73'
Set-PlainValue 'E39' '  -1.19%  '
Set-TextValue 'D40' '2.75'
Set-PlainValue 'E40' '  -1.26%  '
Set-TextValue 'D41' '1.62'
Set-PlainValue 'E41' '  +11.80%  '
Set-TextValue 'D42' '97.63'
Set-PlainValue 'E42' '  -1.70%  '
Set-TextValue 'D43' '70.43'
Set-PlainValue 'E43' '  +1.60%  '
Set-PlainValue 'E44' '  +0.62%  '
Set-PlainValue 'E45' '  -0.49%  '
Set-TextValue 'D46' '12.44'
Set-PlainValue 'E46' '  +3.84%  '
Set-TextValue 'D47' '117.54'
Set-PlainValue 'E47' '  +5.17%  '
Set-TextValue 'D48' '80.91'
Set-PlainValue 'E48' '  +5.25%  '
Set-PlainValue 'D49' '1.643.28'
Set-PlainValue 'E49' '  +3.22%  '
Set-TextValue 'D50' '8.89'
Set-PlainValue 'E50' '  -0.33%  '
Set-PlainValue 'E51' '  -0.03%  '
